$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 3774.1853
$ws.Cells.Item(51, 9).Value = 4833.25
$ws.Cells.Item(51, 10).Value = 2926.9333
$ws.Cells.Item(51, 11).Value = 4833.25
$ws.Cells.Item(51, 12).Value = 2926.9333
$ws.Cells.Item(51, 13).Value = -4349.25
$ws.Cells.Item(51, 14).Value = -3894.9333

$ws.Cells.Item(100, 8).Value = 2375.238
$ws.Cells.Item(100, 9).Value = 1364.5834
$ws.Cells.Item(100, 11).Value = 1364.5834
$ws.Cells.Item(100, 13).Value = -823.5834

$ws.Cells.Item(134, 8).Value = 105000
$ws.Cells.Item(134, 10).Value = 105000
$ws.Cells.Item(134, 12).Value = 105000
$ws.Cells.Item(134, 14).Value = -115140

$ws.Cells.Item(137, 8).Value = 36906.918
$ws.Cells.Item(137, 9).Value = 40932.24
$ws.Cells.Item(137, 11).Value = 122796.72
$ws.Cells.Item(137, 13).Value = -120246.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9610.352999999999
$ws.Cells.Item(32, 9).Value = 5225.067
$ws.Cells.Item(32, 11).Value = 5225.067
$ws.Cells.Item(32, 13).Value = -4938.067

$ws.Cells.Item(45, 8).Value = 5283.1665
$ws.Cells.Item(45, 9).Value = 4857
$ws.Cells.Item(45, 10).Value = 5879.8
$ws.Cells.Item(45, 11).Value = 4857
$ws.Cells.Item(45, 12).Value = 5879.8
$ws.Cells.Item(45, 13).Value = -4480
$ws.Cells.Item(45, 14).Value = -6633.8

$ws.Cells.Item(61, 8).Value = 3133.9333
$ws.Cells.Item(61, 9).Value = 2964.9285
$ws.Cells.Item(61, 11).Value = 2964.9285
$ws.Cells.Item(61, 13).Value = -2752.9285

$ws.Cells.Item(110, 8).Value = 1694.5625
$ws.Cells.Item(110, 9).Value = 1185.6154
$ws.Cells.Item(110, 11).Value = 1185.6154
$ws.Cells.Item(110, 13).Value = 859.3846000000001

$ws.Cells.Item(122, 8).Value = 4859.3335
$ws.Cells.Item(122, 9).Value = 2926.6924
$ws.Cells.Item(122, 11).Value = 8780.0772
$ws.Cells.Item(122, 13).Value = -6330.0772

$ws.Cells.Item(132, 8).Value = 296853.47
$ws.Cells.Item(132, 9).Value = 347294.6
$ws.Cells.Item(132, 11).Value = 1041883.8
$ws.Cells.Item(132, 13).Value = -1039353.8

$ws.Cells.Item(136, 8).Value = 3133.9333
$ws.Cells.Item(136, 9).Value = 2964.9285
$ws.Cells.Item(136, 11).Value = 8894.7855
$ws.Cells.Item(136, 13).Value = -6344.7855

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 4020000
$ws.Cells.Item(7, 9).Value = 6676666.5
$ws.Cells.Item(7, 10).Value = 35000
$ws.Cells.Item(7, 11).Value = 6676666.5
$ws.Cells.Item(7, 12).Value = 35000
$ws.Cells.Item(7, 13).Value = -6676553.5
$ws.Cells.Item(7, 14).Value = -35226

$ws.Cells.Item(94, 8).Value = 2291.0715
$ws.Cells.Item(94, 9).Value = 775.1111
$ws.Cells.Item(94, 11).Value = 775.1111
$ws.Cells.Item(94, 13).Value = -324.1111

$ws.Cells.Item(99, 8).Value = 3249.8333
$ws.Cells.Item(99, 9).Value = 3062.25
$ws.Cells.Item(99, 11).Value = 3062.25
$ws.Cells.Item(99, 13).Value = -1564.25

$ws.Cells.Item(107, 8).Value = 2543.6428
$ws.Cells.Item(107, 9).Value = 1451.375
$ws.Cells.Item(107, 10).Value = 4000
$ws.Cells.Item(107, 11).Value = 1451.375
$ws.Cells.Item(107, 12).Value = 4000
$ws.Cells.Item(107, 13).Value = 468.625
$ws.Cells.Item(107, 14).Value = -7840

$ws.Cells.Item(134, 8).Value = 9806026
$ws.Cells.Item(134, 9).Value = 1662.3572
$ws.Cells.Item(134, 11).Value = 4987.071599999999
$ws.Cells.Item(134, 13).Value = -2452.071599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2993.675
$ws.Cells.Item(31, 9).Value = 2152.5
$ws.Cells.Item(31, 10).Value = 4255.4375
$ws.Cells.Item(31, 11).Value = 2152.5
$ws.Cells.Item(31, 12).Value = 4255.4375
$ws.Cells.Item(31, 13).Value = -1857.5
$ws.Cells.Item(31, 14).Value = -4845.4375

$ws.Cells.Item(34, 8).Value = 2993.675
$ws.Cells.Item(34, 9).Value = 2152.5
$ws.Cells.Item(34, 10).Value = 4255.4375
$ws.Cells.Item(34, 11).Value = 2152.5
$ws.Cells.Item(34, 12).Value = 4255.4375
$ws.Cells.Item(34, 13).Value = -1950.5
$ws.Cells.Item(34, 14).Value = -4659.4375

$ws.Cells.Item(41, 8).Value = 3251
$ws.Cells.Item(41, 9).Value = 3251
$ws.Cells.Item(41, 11).Value = 3251
$ws.Cells.Item(41, 13).Value = -2823

$ws.Cells.Item(58, 8).Value = 3280.389
$ws.Cells.Item(58, 9).Value = 2964.7827
$ws.Cells.Item(58, 11).Value = 2964.7827
$ws.Cells.Item(58, 13).Value = -2761.7827

$ws.Cells.Item(68, 8).Value = 65500
$ws.Cells.Item(68, 10).Value = 65500
$ws.Cells.Item(68, 12).Value = 65500
$ws.Cells.Item(68, 14).Value = -66998

$ws.Cells.Item(71, 8).Value = 65500
$ws.Cells.Item(71, 10).Value = 65500
$ws.Cells.Item(71, 12).Value = 196500
$ws.Cells.Item(71, 14).Value = -203988

$ws.Cells.Item(86, 8).Value = 6603.5
$ws.Cells.Item(86, 9).Value = 9007
$ws.Cells.Item(86, 10).Value = 4200
$ws.Cells.Item(86, 11).Value = 9007
$ws.Cells.Item(86, 12).Value = 4200
$ws.Cells.Item(86, 13).Value = -7884
$ws.Cells.Item(86, 14).Value = -6446

$ws.Cells.Item(89, 8).Value = 6603.5
$ws.Cells.Item(89, 9).Value = 9007
$ws.Cells.Item(89, 10).Value = 4200
$ws.Cells.Item(89, 11).Value = 45035
$ws.Cells.Item(89, 12).Value = 21000
$ws.Cells.Item(89, 13).Value = -39419
$ws.Cells.Item(89, 14).Value = -32232

$ws.Cells.Item(99, 8).Value = 4122.25
$ws.Cells.Item(99, 9).Value = 3829.6667
$ws.Cells.Item(99, 11).Value = 3829.6667
$ws.Cells.Item(99, 13).Value = -2331.6667

$ws.Cells.Item(107, 8).Value = 79087.38
$ws.Cells.Item(107, 9).Value = 126154.5
$ws.Cells.Item(107, 11).Value = 126154.5
$ws.Cells.Item(107, 13).Value = -124234.5

$ws.Cells.Item(126, 8).Value = 4122.25
$ws.Cells.Item(126, 9).Value = 3829.6667
$ws.Cells.Item(126, 11).Value = 11489.0001
$ws.Cells.Item(126, 13).Value = -9019.000100000001

$ws.Cells.Item(134, 8).Value = 2117.7856
$ws.Cells.Item(134, 9).Value = 1521.6177
$ws.Cells.Item(134, 10).Value = 4651.5
$ws.Cells.Item(134, 11).Value = 4564.8531
$ws.Cells.Item(134, 12).Value = 13954.5
$ws.Cells.Item(134, 13).Value = -2029.8531
$ws.Cells.Item(134, 14).Value = -19024.5

$ws.Cells.Item(136, 8).Value = 3280.389
$ws.Cells.Item(136, 9).Value = 2964.7827
$ws.Cells.Item(136, 11).Value = 8894.348100000001
$ws.Cells.Item(136, 13).Value = -6344.348100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 17
$ws.Cells.Item(12, 10).Value = 14
$ws.Cells.Item(12, 12).Value = 42
$ws.Cells.Item(12, 14).Value = -388

$ws.Cells.Item(110, 8).Value = 7
$ws.Cells.Item(110, 9).Value = 7
$ws.Cells.Item(110, 11).Value = 21
$ws.Cells.Item(110, 13).Value = 4069

$ws.Cells.Item(132, 8).Value = 946.8333
$ws.Cells.Item(132, 9).Value = 916.4
$ws.Cells.Item(132, 11).Value = 8247.6
$ws.Cells.Item(132, 13).Value = -5717.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 652.375
$ws.Cells.Item(107, 10).Value = 500
$ws.Cells.Item(107, 12).Value = 500
$ws.Cells.Item(107, 14).Value = -4340

$ws.Cells.Item(132, 8).Value = 3229.4285
$ws.Cells.Item(132, 9).Value = 3229.4285
$ws.Cells.Item(132, 11).Value = 9688.2855
$ws.Cells.Item(132, 13).Value = -7158.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2500.4443
$ws.Cells.Item(46, 10).Value = 3000.5715
$ws.Cells.Item(46, 12).Value = 3000.5715
$ws.Cells.Item(46, 14).Value = -3376.5715

$ws.Cells.Item(61, 8).Value = 3454.3333
$ws.Cells.Item(61, 9).Value = 1198.2222
$ws.Cells.Item(61, 10).Value = 5710.4443
$ws.Cells.Item(61, 11).Value = 1198.2222
$ws.Cells.Item(61, 12).Value = 5710.4443
$ws.Cells.Item(61, 13).Value = -996.2221999999999
$ws.Cells.Item(61, 14).Value = -6114.4443

$ws.Cells.Item(93, 8).Value = 1749.625
$ws.Cells.Item(93, 9).Value = 1109
$ws.Cells.Item(93, 10).Value = 1963.1666
$ws.Cells.Item(93, 11).Value = 1109
$ws.Cells.Item(93, 12).Value = 1963.1666
$ws.Cells.Item(93, 13).Value = 139
$ws.Cells.Item(93, 14).Value = -4459.1666

$ws.Cells.Item(113, 8).Value = 3454.3333
$ws.Cells.Item(113, 9).Value = 1198.2222
$ws.Cells.Item(113, 10).Value = 5710.4443
$ws.Cells.Item(113, 11).Value = 1198.2222
$ws.Cells.Item(113, 12).Value = 5710.4443
$ws.Cells.Item(113, 13).Value = 971.7778000000001
$ws.Cells.Item(113, 14).Value = -10050.4443

$ws.Cells.Item(132, 8).Value = 387440.47
$ws.Cells.Item(132, 10).Value = 2732.3333
$ws.Cells.Item(132, 12).Value = 8196.999899999999
$ws.Cells.Item(132, 14).Value = -13256.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 370.9091
$ws.Cells.Item(107, 9).Value = 276.9091
$ws.Cells.Item(107, 11).Value = 830.7273
$ws.Cells.Item(107, 13).Value = 1089.2727

$ws.Cells.Item(113, 8).Value = 274.375
$ws.Cells.Item(113, 9).Value = 341
$ws.Cells.Item(113, 10).Value = 163.33333
$ws.Cells.Item(113, 11).Value = 1023
$ws.Cells.Item(113, 12).Value = 489.99999
$ws.Cells.Item(113, 13).Value = 1147
$ws.Cells.Item(113, 14).Value = -4829.99999

$ws.Cells.Item(132, 8).Value = 27157.232
$ws.Cells.Item(132, 9).Value = 32948.883
$ws.Cells.Item(132, 11).Value = 98846.649
$ws.Cells.Item(132, 13).Value = -96316.649
